$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 92 (pushes existing rows 92-120 down to 93-121)
$ws.Rows(92).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A92").Value = 10
$ws.Range("B92").Value = "Vega Modelo de Temuco"
$ws.Range("C92").Value = "La Araucanía"
$ws.Range("D92").Value = 45120
$ws.Range("E92").Value = 9
$ws.Range("F92").Value = 300000001
$ws.Range("G92").Value = "Rabanito"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 80
$ws.Range("K92").Value = 7000
$ws.Range("L92").Value = 7000
$ws.Range("M92").Value = 7000
$ws.Range("N92").Value = "`$/docena de paquetes"
$ws.Range("O92").Value = "Provincia de Cautín"
$ws.Range("P92").Value = 583
$ws.Range("Q92").Value = 12
$ws.Range("R92").Value = "Hortaliza"
